$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NFL")
Write-Host $ws.Name
